$wb = $excel.ActiveWorkbook

# --- About sheet: add Minnesota label + date stamp next to the title ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("B1").Value = "Minnesota"
$wsAbout.Range("C1").Value = (Get-Date -Year 2021 -Month 9 -Day 24 -Hour 0 -Minute 0 -Second 0).Date
$wsAbout.Range("C1").NumberFormat = "m/d/yyyy"

# --- BDPbES sheet: update priority order values for MN (1 -> 2) ---
$wsData = $wb.Worksheets.Item("BDPbES")
$wsData.Range("B5").Value = 2   # hydro
$wsData.Range("B8").Value = 2   # solar thermal
$wsData.Range("B9").Value = 2   # biomass
$wsData.Range("B10").Value = 2  # geothermal
$wsData.Range("B14").Value = 2  # offshore wind

# keep B8 as the active selection on the data sheet, matching the saved view
$wsData.Range("B8").Select()
